# Insert a new data row at row 73 (shifting the existing rows 73-117 down
# to 74-118) and populate it with the new record, matching the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 73..117 down by one row.
$ws.Rows.Item(73).Insert()

# Populate the newly inserted row 73 with the new observation.
$ws.Cells.Item(73, 1).Value  = 9
$ws.Cells.Item(73, 2).Value  = 'Vega Central Mapocho de Santiago'
$ws.Cells.Item(73, 3).Value  = 'Metropolitana'
$ws.Cells.Item(73, 4).Value  = 44767
$ws.Cells.Item(73, 5).Value  = 13
$ws.Cells.Item(73, 6).Value  = 100112022
$ws.Cells.Item(73, 7).Value  = 'Arveja Verde'
$ws.Cells.Item(73, 8).Value  = 'Perfection'
$ws.Cells.Item(73, 9).Value  = 'Primera'
$ws.Cells.Item(73, 10).Value = 16
$ws.Cells.Item(73, 11).Value = 40000
$ws.Cells.Item(73, 12).Value = 40000
$ws.Cells.Item(73, 13).Value = 40000
$ws.Cells.Item(73, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(73, 15).Value = 'Provincia de Huasco'
$ws.Cells.Item(73, 16).Value = 1600
$ws.Cells.Item(73, 17).Value = 25
$ws.Cells.Item(73, 18).Value = 'Hortaliza'
